$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.41499390233946015227
$ws.Range("A2").Value = -0.22966111487525239787
$ws.Range("A3").Value = -0.22454460669028558639
$ws.Range("A4").Value = 0.00016154075696733601
$ws.Range("A5").Value = 0.00024984267888622578
$ws.Range("A6").Value = -0.00039288925875659964
$ws.Range("A7").Value = -0.00023740224794665870
$ws.Range("A8").Value = -0.00011731229479506278
$ws.Range("A9").Value = 0.00021656415369121303
$ws.Range("A10").Value = -0.00011074137397429396
$ws.Range("A11").Value = -0.00040554697075647570
$ws.Range("A12").Value = -0.00022532072707057924
$ws.Range("A13").Value = -0.00000911761144039945
$ws.Range("A14").Value = -0.00029767619429069413
$ws.Range("A15").Value = -0.00080934271602350091
$ws.Range("A16").Value = 0.00079199936049564965
$ws.Range("A17").Value = 0.00088533466953085685
$ws.Range("A18").Value = 0.00068761074348950419
$ws.Range("A19").Value = 0.00080940944983294376
$ws.Range("A20").Value = 0.00050470023674740127
$ws.Range("A21").Value = 0.00071838967533614459
$ws.Range("A22").Value = 0.00032059039849497399
$ws.Range("A23").Value = 0.00123990437220014007
$ws.Range("A24").Value = -0.00028003090355892588
